$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap Pin_Sequencer_SCL (J9) and Pin_Sequencer_SDA (J10)
$ws.Range("J9").Value = "Pin_Sequencer_SDA"
$ws.Range("J10").Value = "Pin_Sequencer_SCL"

# Remove Pin_SW3_in (E10) and Pin_SW4_in (E11) pin assignments
$ws.Range("E10").ClearContents()
$ws.Range("E11").ClearContents()

# Swap LED_RED (J22) and LED_GREEN (J24)
$ws.Range("J22").Value = "LED_GREEN"
$ws.Range("J24").Value = "LED_RED"

# Update the active selection left by the editor
$ws.Range("J20").Select()
